$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "word"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "is_active"

$data = @(
    @("eng", "shit",   "Blacklisted Word", $true),
    @("eng", "damn",   "Blacklisted Word", $true),
    @("eng", "nigga",  "Blacklisted Word", $true),
    @("eng", "dammit", "Blacklisted Word", $true),
    @("fra", "Merde",  "Mot sur la liste noire", $true),
    @("fra", "Damn",   "Mot sur la liste noire", $true),
    @("fra", "nigga",  "Mot sur la liste noire", $true),
    @("fra", "bon sang", "Mot sur la liste noire", $true)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Apply the header (column A) style from A1 to A2:A9, matching s="1" on those cells.
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
